# Append the 2025-03-06 price row to each of the Argent (silver) price
# sheets, mirroring the existing Date/Price text-cell layout.
$wb = $excel.ActiveWorkbook

$newDate = "2025-03-06"

# Sheet name -> new Price value for row 5 (all stored as text, like the
# existing rows in these sheets).
$updates = [ordered]@{
    "N-type Wafer"             = "1.19"
    "Cell Topcon 183mm"        = "0.293"
    "Module Topcon 183mm"      = "0.1"
    "Silver Rear_side"         = "5,211"
    "Silver Busbar front-side" = "7,801"
    "Silver finger front-side" = "7,851"
    "USD_CNY"                  = "7.2842"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    $dateCell = $ws.Range("A5")
    $priceCell = $ws.Range("B5")

    # Force text storage (matching the other rows' text cells) instead of
    # letting Excel auto-detect a date/number type.
    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"

    $dateCell.Value = $newDate
    $priceCell.Value = $updates[$sheetName]
}
